$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6 - shifts existing rows 6..41 down to 7..42
$ws.Rows.Item(6).Insert()

# Match formatting of the rank column by copying the style from the row above
$ws.Range("A5").Copy($ws.Range("A6"))

# Populate the newly inserted row with the new researcher's data
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Jawhar"
$ws.Range("C6").Value = "Hafsa"
$ws.Range("D6").Value = "Mohammed VI Polytechnic University"
$ws.Range("E6").Value = "Morocco"
$ws.Range("F6").Value = "ww5A_WMAAAAJ"
$ws.Range("G6").Value = "M"
$ws.Range("H6").Value = 1987
$ws.Range("I6").Value = "Médecine, Biologie et Sciences de la Santé"
$ws.Range("J6").Value = 17
$ws.Range("K6").Value = 21
$ws.Range("L6").Value = 16
$ws.Range("M6").Value = 21
$ws.Range("N6").Value = 1159
$ws.Range("O6").Value = 1069

# Renumber the rank column (A) for all rows pushed down by the insert
for ($r = 7; $r -le 42; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

Write-Output "done"
